$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update Puntaje (C) values for summary rows
$ws.Range("C2").Value = 0
$ws.Range("C3").Value = 0.2
$ws.Range("C4").Value = 1
$ws.Range("C5").Value = 0.8

# Update Observaciones (D) for the first two summary rows
$ws.Range("D3").Value = "Se puede simplificar los If"
$ws.Range("D4").Value = "Falla en:qwertyuiop, 01012001"

# Update Cumple (B) and Observaciones (D) for detail rows that changed from Si to No
# (entered bottom-up, as in the original edit)
$ws.Range("B13").Value = "No"

$ws.Range("D24").Value = "Demasiados If y condiciones largas, en lugar de usar un Switch"
$ws.Range("B24").Value = "No"

$ws.Range("D23").Value = "No hqay retornos directos, hay if and else"
$ws.Range("B23").Value = "No"

$ws.Range("D22").Value = "Oeracionees de Validacion combinadas"
$ws.Range("B22").Value = "No"

$ws.Range("D21").Value = "Hay mucha logica combinada"
$ws.Range("B21").Value = "No"

$ws.Range("D13").Value = "No hay uso de constantse donde si pudiera tnerlas"

# Last entered: Observaciones for row 5
$ws.Range("D5").Value = "Se complica en algunos momentos por la combinacion de logica en las funciones"

# Update the active selection cell
$ws.Range("D6").Select()
